$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 0.0003472222222222222
$ws.Range("K2").Value = 2487
$ws.Range("L2").Value = 0.004974
